# Generate Report for Handback
# Update the "handed back" / xliff-generation timestamps that were
# refreshed by the latest report run.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the in-sync-with-en-US row
$wsOverview.Range("G3").Value = "2016-08-26 18:47:42"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime
$wsZhCn.Range("H3").Value = "2016-08-26 18:47:37"
$wsZhCn.Range("K3").Value = "2016-08-26 18:47:55"

# de-de sheet: Correspond Handoff Datetime / Correspond Handback DateTime
$wsDeDe.Range("H3").Value = "2016-08-26 18:47:42"
$wsDeDe.Range("K3").Value = "2016-08-26 18:48:09"
